# cryptos.xlsx refresh — Mon Dec 11 18:09:20 UTC 2023 (GitHub Actions)
#
# Updates the per-coin Price (D) and Volume(1h) (E) figures pulled from
# coinranking.com, and re-syncs a handful of rows whose rank (and thus
# Coin/Link/Price/Volume) shifted in the source feed (rows 43/44, 46/47, 51).
#
# All D/E/B/C cells in this sheet are stored as plain text (inlineStr) in the
# source workbook, not numbers — the "%"-suffixed volume strings and the
# dotted/odd-format prices (e.g. "41.267.12", "0.0₃0942") are literal text.
# Most new price strings (e.g. "243.00", "1.00") *look* numeric, though, and
# a bare `.Value = "243.00"` assignment would make Excel auto-coerce them to
# the number 243 (dropping the trailing zero) and flip the cell to numeric
# storage. To keep them as exact text we briefly mark those cells as Text
# (NumberFormat "@") before writing, then restore the default "Normal" style
# afterwards so no stray number-format style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / percentage / coin-name / link updates ---
# (never numeric-looking, so a direct .Value assignment keeps them as text)
$ws.Range('D2').Value = '41.267.12'
$ws.Range('E2').Value = '  -5.88%  '
$ws.Range('D3').Value = '2.198.28'
$ws.Range('E3').Value = '  -6.42%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('E6').Value = '  -6.77%  '
$ws.Range('E7').Value = '  -3.84%  '
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  -9.92%  '
$ws.Range('E10').Value = '  +9.73%  '
$ws.Range('E11').Value = '  -6.78%  '
$ws.Range('E12').Value = '  -5.16%  '
$ws.Range('E13').Value = '  -3.51%  '
$ws.Range('E14').Value = '  -8.59%  '
$ws.Range('D15').Value = '2.529.42'
$ws.Range('E15').Value = '  -6.16%  '
$ws.Range('E16').Value = '  -9.35%  '
$ws.Range('E17').Value = '  -7.75%  '
$ws.Range('D18').Value = '2.204.63'
$ws.Range('E18').Value = '  -5.99%  '
$ws.Range('D19').Value = '41.272.69'
$ws.Range('E19').Value = '  -5.78%  '
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  -8.34%  '
$ws.Range('E21').Value = '  -5.02%  '
$ws.Range('E22').Value = '  -7.08%  '
$ws.Range('E23').Value = '  -7.80%  '
$ws.Range('E24').Value = '  +8.81%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('E26').Value = '  -4.98%  '
$ws.Range('E27').Value = '  -2.89%  '
$ws.Range('E28').Value = '  -3.92%  '
$ws.Range('E29').Value = '  -6.81%  '
$ws.Range('E30').Value = '  -4.23%  '
$ws.Range('E31').Value = '  -8.85%  '
$ws.Range('E32').Value = '  -7.13%  '
$ws.Range('E33').Value = '  -7.21%  '
$ws.Range('E34').Value = '  -5.13%  '
$ws.Range('E35').Value = '  -4.54%  '
$ws.Range('E36').Value = '  -9.24%  '
$ws.Range('E37').Value = '  +3.62%  '
$ws.Range('E38').Value = '  +17.89%  '
$ws.Range('E39').Value = '  -6.21%  '
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('E41').Value = '  -10.01%  '
$ws.Range('E42').Value = '  -1.32%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E43').Value = '  -3.37%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E44').Value = '  -14.47%  '
$ws.Range('E45').Value = '  -2.97%  '
$ws.Range('B46').Value = 'BinanceUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E47').Value = '  -6.59%  '
$ws.Range('E48').Value = '  +5.07%  '
$ws.Range('E49').Value = '  +5.21%  '
$ws.Range('E50').Value = '  -5.34%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E51').Value = '  -5.58%  '

# --- Numeric-looking Price (D) updates ---
# Force text storage so values like "243.00" / "1.00" keep their exact
# formatting instead of being coerced into numbers 243 / 1.
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '243.00'
$cell.Style = "Normal"
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '0.620'
$cell.Style = "Normal"
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '70.16'
$cell.Style = "Normal"
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.540'
$cell.Style = "Normal"
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '36.70'
$cell.Style = "Normal"
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.0944'
$cell.Style = "Normal"
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '57.75'
$cell.Style = "Normal"
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '6.60'
$cell.Style = "Normal"
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '14.63'
$cell.Style = "Normal"
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '0.834'
$cell.Style = "Normal"
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '73.75'
$cell.Style = "Normal"
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '6.02'
$cell.Style = "Normal"
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '232.52'
$cell.Style = "Normal"
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '2.41'
$cell.Style = "Normal"
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '9.71'
$cell.Style = "Normal"
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '168.38'
$cell.Style = "Normal"
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '20.21'
$cell.Style = "Normal"
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '0.118'
$cell.Style = "Normal"
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '0.0702'
$cell.Style = "Normal"
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '4.57'
$cell.Style = "Normal"
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '3.89'
$cell.Style = "Normal"
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '23.25'
$cell.Style = "Normal"
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '2.25'
$cell.Style = "Normal"
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '0.0269'
$cell.Style = "Normal"
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '5.79'
$cell.Style = "Normal"
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '64.51'
$cell.Style = "Normal"
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '8.76'
$cell.Style = "Normal"
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '4.77'
$cell.Style = "Normal"
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '0.0988'
$cell.Style = "Normal"
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '4.49'
$cell.Style = "Normal"
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '10.02'
$cell.Style = "Normal"
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '1.16'
$cell.Style = "Normal"
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '1.08'
$cell.Style = "Normal"
